# Applies the "Voeg kennis van KA toe aan M23" edit set:
#  1. Slide 13 ("M16: Het project gebruikt tools voor vastgestelde taken"):
#     remove the bullet about "controleren van de configuratie op aanwezigheid
#     van bekende kwetsbaarheden," (OpenVAS reference removal).
#  2. Slide 19 ("M23: ..."):
#     - title gains "kennis van en" before "ervaring met de Kwaliteitsaanpak"
#     - body paragraph gains an extra sentence about onboarding new project
#       members on the Kwaliteitsaanpak.

$p = $ppt.ActivePresentation

# --- 1. Slide 13: drop the "controleren van de configuratie ..." bullet ---
$s13 = $p.Slides.Item(13)
$bulletsShape13 = $s13.Shapes.Item("TextBox 2")
$tr13 = $bulletsShape13.TextFrame2.TextRange
$fullText13 = $tr13.Text
$target13 = "controleren van de configuratie op aanwezigheid van bekende kwetsbaarheden,"
$idx13 = $fullText13.IndexOf($target13)
if ($idx13 -ge 0) {
    # +1 char to also swallow the trailing paragraph-break so the whole
    # paragraph (not just its text) disappears, keeping the other bullets'
    # <a:pPr> formatting untouched.
    $sub13 = $tr13.Characters($idx13 + 1, $target13.Length + 1)
    $sub13.Delete()
}

# --- 2. Slide 19: title + body updates ---
$s19 = $p.Slides.Item(19)

$titleShape19 = $s19.Shapes.Item("Title 1")
$titleShape19.TextFrame.TextRange.Text = "M23: Het project zorgt voor de aanwezigheid van kennis van en ervaring met de Kwaliteitsaanpak"

$bodyShape19 = $s19.Shapes.Item("TextBox 2")
$bodyShape19.TextFrame.TextRange.Text = "De software delivery manager zorgt ervoor dat bij nieuwe projecten wordt gestart met ten minste twee projectleden die bekend zijn met de Kwaliteitsaanpak. Projectleden die nog niet bekend zijn met de Kwaliteitsaanpak krijgen uitleg over de inhoud en achtergrond van de Kwaliteitsaanpak."
